$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 551
$ws.Range("I12").Value = 609.5
$ws.Range("K12").Value = 609.5
$ws.Range("M12").Value = -439.5
$ws.Range("H39").Value = 220.42105
$ws.Range("I39").Value = 115.52941
$ws.Range("J39").Value = 1112
$ws.Range("K39").Value = 346.58823
$ws.Range("L39").Value = 3336
$ws.Range("M39").Value = -50.58823000000001
$ws.Range("N39").Value = -3928
$ws.Range("H100").Value = 1846.4375
$ws.Range("I100").Value = 1413.091
$ws.Range("K100").Value = 1413.091
$ws.Range("M100").Value = -872.0909999999999
$ws.Range("H116").Value = 28411.277
$ws.Range("I116").Value = 18958.166
$ws.Range("K116").Value = 18958.166
$ws.Range("M116").Value = -15516.166
$ws.Range("H127").Value = 80248.69500000001
$ws.Range("I127").Value = 94203
$ws.Range("J127").Value = 3500
$ws.Range("K127").Value = 282609
$ws.Range("L127").Value = 10500
$ws.Range("M127").Value = -277649
$ws.Range("N127").Value = -20420
$ws.Range("H132").Value = 6990.436
$ws.Range("I132").Value = 7327.757
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 21983.271
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -19453.271
$ws.Range("N132").Value = -7310
$ws.Range("H137").Value = 20005860
$ws.Range("I137").Value = 26317686
$ws.Range("J137").Value = 18416.334
$ws.Range("K137").Value = 78953058
$ws.Range("L137").Value = 55249.00199999999
$ws.Range("M137").Value = -78950508
$ws.Range("N137").Value = -60349.00199999999
$ws.Range("H138").Value = 8387.321
$ws.Range("I138").Value = 8947.380999999999
$ws.Range("K138").Value = 26842.143
$ws.Range("M138").Value = -21702.143

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 706.7692
$ws.Range("I2").Value = 640.6667
$ws.Range("K2").Value = 640.6667
$ws.Range("M2").Value = -527.6667
$ws.Range("H32").Value = 696217.4399999999
$ws.Range("I32").Value = 751998.4399999999
$ws.Range("K32").Value = 751998.4399999999
$ws.Range("M32").Value = -751711.4399999999
$ws.Range("H47").Value = 6000
$ws.Range("I47").Value = 6000
$ws.Range("K47").Value = 6000
$ws.Range("M47").Value = -5275
$ws.Range("H116").Value = 706.7692
$ws.Range("I116").Value = 640.6667
$ws.Range("K116").Value = 640.6667
$ws.Range("M116").Value = 1653.3333
$ws.Range("H132").Value = 864797.75
$ws.Range("I132").Value = 964074.4399999999
$ws.Range("J132").Value = 4399.6665
$ws.Range("K132").Value = 2892223.32
$ws.Range("L132").Value = 13198.9995
$ws.Range("M132").Value = -2889693.32
$ws.Range("N132").Value = -18258.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 706.7692
$ws.Range("I3").Value = 640.6667
$ws.Range("K3").Value = 640.6667
$ws.Range("M3").Value = -526.6667
$ws.Range("H105").Value = 2356.0833
$ws.Range("I105").Value = 1153.5714
$ws.Range("J105").Value = 4039.6
$ws.Range("K105").Value = 1153.5714
$ws.Range("L105").Value = 4039.6
$ws.Range("M105").Value = 593.4286
$ws.Range("N105").Value = -7533.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1534.3334
$ws.Range("I22").Value = 225.22223
$ws.Range("J22").Value = 2516.1667
$ws.Range("K22").Value = 225.22223
$ws.Range("L22").Value = 2516.1667
$ws.Range("M22").Value = 124.77777
$ws.Range("N22").Value = -3216.1667
$ws.Range("H58").Value = 6638243.5
$ws.Range("I58").Value = 5955893
$ws.Range("K58").Value = 5955893
$ws.Range("M58").Value = -5955690
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 1288.7
$ws.Range("I132").Value = 1288.7
$ws.Range("K132").Value = 3866.1
$ws.Range("M132").Value = -1336.1
$ws.Range("H134").Value = 1636.9028
$ws.Range("I134").Value = 1375.0741
$ws.Range("J134").Value = 2422.389
$ws.Range("K134").Value = 4125.2223
$ws.Range("L134").Value = 7267.167
$ws.Range("M134").Value = -1590.2223
$ws.Range("N134").Value = -12337.167
$ws.Range("H136").Value = 6638243.5
$ws.Range("I136").Value = 5955893
$ws.Range("K136").Value = 17867679
$ws.Range("M136").Value = -17865129

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2269409.5
$ws.Range("J5").Value = 1939682.9
$ws.Range("L5").Value = 5819048.699999999
$ws.Range("N5").Value = -5819272.699999999
$ws.Range("H17").Value = 375
$ws.Range("J17").Value = 485
$ws.Range("L17").Value = 1455
$ws.Range("N17").Value = -1793
$ws.Range("H19").Value = 129.16667
$ws.Range("I19").Value = 129.16667
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 387.50001
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -213.50001
$ws.Range("N19").ClearContents()
$ws.Range("H25").Value = 324.5
$ws.Range("I25").Value = 324.5
$ws.Range("K25").Value = 973.5
$ws.Range("M25").Value = -804.5
$ws.Range("H29").Value = 441.8
$ws.Range("J29").Value = 25
$ws.Range("L29").Value = 75
$ws.Range("N29").Value = -629
$ws.Range("H30").Value = 324.5
$ws.Range("I30").Value = 324.5
$ws.Range("K30").Value = 973.5
$ws.Range("M30").Value = -871.5
$ws.Range("H31").Value = 5950.5
$ws.Range("J31").Value = 1900
$ws.Range("L31").Value = 5700
$ws.Range("N31").Value = -6276
$ws.Range("H70").Value = 2198.3333
$ws.Range("I70").Value = 2297.5
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 6892.5
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -6577.5
$ws.Range("N70").Value = -6630
$ws.Range("H73").Value = 2198.3333
$ws.Range("I73").Value = 2297.5
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 6892.5
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -5800.5
$ws.Range("N73").Value = -8184
$ws.Range("H135").Value = 2269409.5
$ws.Range("J135").Value = 1939682.9
$ws.Range("L135").Value = 17457146.1
$ws.Range("N135").Value = -17462216.1
$ws.Range("H138").Value = 9595.5
$ws.Range("I138").Value = 9494.625
$ws.Range("J138").Value = 9999
$ws.Range("K138").Value = 28483.875
$ws.Range("L138").Value = 29997
$ws.Range("M138").Value = -23343.875
$ws.Range("N138").Value = -40277

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.75
$ws.Range("I2").Value = 59.75
$ws.Range("K2").Value = 59.75
$ws.Range("M2").Value = 53.25
$ws.Range("H80").Value = 2359.25
$ws.Range("I80").Value = 2216.2856
$ws.Range("J80").Value = 2692.8333
$ws.Range("K80").Value = 2216.2856
$ws.Range("L80").Value = 2692.8333
$ws.Range("M80").Value = -1218.2856
$ws.Range("N80").Value = -4688.8333
$ws.Range("H83").Value = 2359.25
$ws.Range("I83").Value = 2216.2856
$ws.Range("J83").Value = 2692.8333
$ws.Range("K83").Value = 11081.428
$ws.Range("L83").Value = 13464.1665
$ws.Range("M83").Value = -6089.428
$ws.Range("N83").Value = -23448.1665
$ws.Range("H126").Value = 7011.875
$ws.Range("I126").Value = 8082.5
$ws.Range("K126").Value = 24247.5
$ws.Range("M126").Value = -21777.5
$ws.Range("H132").Value = 16517.088
$ws.Range("J132").Value = 1007.5
$ws.Range("L132").Value = 3022.5
$ws.Range("N132").Value = -8082.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5356.2144
$ws.Range("I61").Value = 5493.304
$ws.Range("K61").Value = 5493.304
$ws.Range("M61").Value = -5291.304
$ws.Range("H113").Value = 5356.2144
$ws.Range("I113").Value = 5493.304
$ws.Range("K113").Value = 5493.304
$ws.Range("M113").Value = -3323.304
$ws.Range("H122").Value = 6878.5
$ws.Range("I122").Value = 4266.6665
$ws.Range("K122").Value = 12799.9995
$ws.Range("M122").Value = -10349.9995
$ws.Range("H133").Value = 89315.336
$ws.Range("J133").Value = 89315.336
$ws.Range("L133").Value = 89315.336
$ws.Range("N133").Value = -94375.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -4002
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 76927330
$ws.Range("I81").Value = 4398.25
$ws.Range("J81").Value = 200004020
$ws.Range("K81").Value = 8796.5
$ws.Range("L81").Value = 400008040
$ws.Range("M81").Value = -7735.5
$ws.Range("N81").Value = -400010162
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 76927330
$ws.Range("I84").Value = 4398.25
$ws.Range("J84").Value = 200004020
$ws.Range("K84").Value = 43982.5
$ws.Range("L84").Value = 2000040200
$ws.Range("M84").Value = -38678.5
$ws.Range("N84").Value = -2000050808
$ws.Range("H113").Value = 651.5789
$ws.Range("I113").Value = 544.2143
$ws.Range("J113").Value = 952.2
$ws.Range("K113").Value = 1632.6429
$ws.Range("L113").Value = 2856.6
$ws.Range("M113").Value = 537.3571000000002
$ws.Range("N113").Value = -7196.6
$ws.Range("H122").Value = 52995.176
$ws.Range("J122").Value = 190911.83
$ws.Range("L122").Value = 572735.49
$ws.Range("N122").Value = -577635.49
$ws.Range("H132").Value = 3146321.8
$ws.Range("I132").Value = 3877230.5
$ws.Range("J132").Value = 3414.8
$ws.Range("K132").Value = 11631691.5
$ws.Range("L132").Value = 10244.4
$ws.Range("M132").Value = -11629161.5
$ws.Range("N132").Value = -15304.4
$ws.Range("H136").Value = 7175950
$ws.Range("I136").Value = 3107337.8
$ws.Range("J136").Value = 16669378
$ws.Range("K136").Value = 9322013.399999999
$ws.Range("L136").Value = 50008134
$ws.Range("M136").Value = -9319463.399999999
$ws.Range("N136").Value = -50013234
